# The "Recorded By" column (G) lists the users who recorded each session as a
# comma-separated string, e.g. "System, dnasr281@gmail.com". This edit
# reverses the order of the names in that list for every data row (row 1 is
# the header), e.g. "dnasr281@gmail.com, System". Rows whose "Recorded By"
# value only contains a single name are left untouched, since reversing a
# single-element list is a no-op.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Count -gt 1) {
        # Manually reverse the list (order: last name first).
        $count = $parts.Count
        $reversedParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newValue = [string]::Join(", ", $reversedParts)
        $cell.Value = $newValue
    }
}
